$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d61a64796dda5162ddf777c7069e4b8c9916e8c5/e2e/2111085c-ee7f-4595-aa5c-a361ee60d92c.md"

# ---------------------------------------------------------------------------
# Overview sheet: append row 3 for the newly handed-off file
# ---------------------------------------------------------------------------
$overview.Range("A3").Value = "2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
$overview.Range("B3").Value = "e2e\2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
$overview.Range("C3").Value = ".md"
$overview.Range("D3").Value = "'"
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 00:45:35"

$overview.Hyperlinks.Add(
    $overview.Range("B3"),
    $newFileUrl,
    [Type]::Missing,
    [Type]::Missing,
    "e2e\2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
) | Out-Null

$overview.ListObjects.Item(1).Resize($overview.Range("A1:G3")) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: append row 3
# ---------------------------------------------------------------------------
$zhcn.Range("A3").Value = "2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "e2e"
$zhcn.Range("E3").Value = "ht"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "2111085c-ee7f-4595-aa5c-a361ee60d92c.b3d58473b2349a2b8be53f257860d2ded59090f1.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-20 00:45:31"
$zhcn.Range("I3").Value = "'"
$zhcn.Range("J3").Value = "'"
$zhcn.Range("K3").Value = "0001-01-01 00:00:00"
$zhcn.Range("L3").Value = "'"
$zhcn.Range("M3").Value = "'True"
$zhcn.Range("N3").Value = "'"
$zhcn.Range("O3").Value = "'False"
$zhcn.Range("P3").Value = "'"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("A3"),
    $newFileUrl,
    [Type]::Missing,
    [Type]::Missing,
    "2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
) | Out-Null

$zhcn.ListObjects.Item(1).Resize($zhcn.Range("A1:P3")) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: append row 3
# ---------------------------------------------------------------------------
$dede.Range("A3").Value = "2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "e2e"
$dede.Range("E3").Value = "ht"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "2111085c-ee7f-4595-aa5c-a361ee60d92c.b3d58473b2349a2b8be53f257860d2ded59090f1.de-de.xlf"
$dede.Range("H3").Value = "2016-08-20 00:45:35"
$dede.Range("I3").Value = "'"
$dede.Range("J3").Value = "'"
$dede.Range("K3").Value = "0001-01-01 00:00:00"
$dede.Range("L3").Value = "'"
$dede.Range("M3").Value = "'True"
$dede.Range("N3").Value = "'"
$dede.Range("O3").Value = "'False"
$dede.Range("P3").Value = "'"

$dede.Hyperlinks.Add(
    $dede.Range("A3"),
    $newFileUrl,
    [Type]::Missing,
    [Type]::Missing,
    "2111085c-ee7f-4595-aa5c-a361ee60d92c.md"
) | Out-Null

$dede.ListObjects.Item(1).Resize($dede.Range("A1:P3")) | Out-Null
